$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update issue log text content ---

# Row 3: "Nulls in country column" issue now quantified (94 nulls) and
# table corrected from geo_lookup to customers
$ws.Range("B3").Value = "94 Nulls in country column"
$ws.Range("C3").Value = "tradeease.customers"

# Rows 8-11: clarify the *_ts issue descriptions with the impact note
$ws.Range("B8").Value  = "purchase_ts is a date and not a timestamp, unable to perform time analysis"
$ws.Range("B9").Value  = "ship_ts is a date and not a timestamp, unable to perform time analysis"
$ws.Range("B10").Value = "delivery_ts is a date and not a timestamp, unable to perform time analysis"
$ws.Range("B11").Value = "refund_ts is a date and not a timestamp, unable to perform time analysis"

# --- Formatting: center-align the "Resolved" column (E) ---
$ws.Range("E2:E11").HorizontalAlignment = -4108  # xlCenter

# --- Selection state ---
$ws.Range("B4").Select()
